# Rename the worksheet from "My Series" to "Data"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Data"

# --- Column D: duplicate column C's header/meta rows (rows 1-10, 12), and add combined series ---
$ws.Range("D1").Value = "Real Estate Investment: Residential: Hebei [SUM(Yes; Yes)+selected(310902301+310902301)]"
$ws.Range("D2").Value = "Test"
$ws.Range("D4").Value = 'Annual, ending "Dec" of each year'
$ws.Range("D5").Value = "RMB mn"
$ws.Range("D6").Value = "NRT SOURCE TEST"
$ws.Range("D7").Value = "Active"
$ws.Range("D8").Value = 310902301
$ws.Range("D9").Value = "SR4825032"
$ws.Range("D12").Value = 35034

# Row 11: label changes from "Function Description" to "Function Information";
# the formula text moves from C11 to D11, and C11 is cleared.
$ws.Range("A11").Value = "Function Information"
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = "SUM(Yes; Yes)+selected(310902301+310902301)"

# Copy styling from column C onto column D for the meta rows (rows 1-12)
$ws.Range("C1:C12").Copy()
$ws.Range("D1:D12").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("D3").Value = ""
$ws.Range("D10").Value = ""

# --- Data rows: existing row 13 shifts down to row 22; 9 new rows inserted above it ---
$ws.Rows("13:21").Insert()
$ws.Range("A13:D21").Font.Bold = $false

$data = @(
  @(35765, 3271.34, 3271.34, 6542.68),
  @(36130, $null, $null, 0),
  @(36495, 6384.48, 6384.48, 12768.96),
  @(36861, 7111.74, 7111.74, 14223.48),
  @(37226, 8354.02, 8354.02, 16708.04),
  @(37591, 10540.67, 10540.67, 21081.34),
  @(37956, 16421.36, 16421.36, 32842.72),
  @(38322, 22352.63, 22352.63, 44705.26),
  @(38687, 29205.22, 29205.22, 58410.44),
  @(39052, 37962.97, 37962.97, 75925.94)
)

$r = 13
foreach ($row in $data) {
    $ws.Range("A$r").Value = $row[0]
    if ($row[1] -ne $null) { $ws.Range("B$r").Value = $row[1] }
    if ($row[2] -ne $null) { $ws.Range("C$r").Value = $row[2] }
    $ws.Range("D$r").Value = $row[3]
    $r = $r + 1
}

# Number formats: date column A, and value columns B:D use "###0.000"
$ws.Range("A13:A22").NumberFormat = "yyyy"
$ws.Range("B13:D22").NumberFormat = "###0.000"
